# Adds 45X for grid batteries; adjusts 45X implementation in transport sector
#
# Content changes captured from the commit:
#  1) FoSfBPPTtC sheet: the 45X phase-down in 2030-2032 (cols K:M) is no
#     longer stepped down (0.375/0.25/0.125) - it now holds at 0.5 like the
#     surrounding years.
#  2) The workbook's active sheet moves from "About" to "FoSfBPPTtC", with
#     the corresponding selections left on each sheet (About -> B6,
#     FoSfBPPTtC -> I20).

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("FoSfBPPTtC")

# --- Data fix: flatten the 2030-2032 phase-down to a constant 0.5 ---------
$wsData.Range("K2").Value = 0.5
$wsData.Range("L2").Value = 0.5
$wsData.Range("M2").Value = 0.5

# --- Selections on each sheet, then make FoSfBPPTtC the active tab --------
[void]$wsAbout.Range("B6").Select()
[void]$wsData.Activate()
[void]$wsData.Range("I20").Select()
